# Update the "family" column (E) for a set of rows from "Circoviridae" to
# "CRESS" to reflect the addition of references for Redondoviridae,
# Smacoviridae, and other CRESS-DNA virus lineages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(14, 15, 16, 17, 18, 19, 20, 21, 23, 24)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "CRESS"
}

# The active selection moved slightly after the edit (still spans the
# same A1:M24 block, just anchored one row lower).
$sel = $excel.Union($ws.Range("D13"), $ws.Range("A1:M24"))
$sel.Select()
